# Correction engine motor + overflow bottles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Brake Master Cylinder (row 5) and Balance bar (row 6) -> "Néant"
$ws.Range("F5").Value = "Néant"
$ws.Range("F6").Value = "Néant"

# Engine assembly (row 8) and Overflow bottles (row 15) -> "En cours Excel"
$ws.Range("F8").Value = "En cours Excel"
$ws.Range("F15").Value = "En cours Excel"

# Update the remembered selection/active cell shown when the sheet was last saved
$ws.Range("G15").Select()
